$d = $word.ActiveDocument

# Locate the "nDaysMinBelowZero_Apr" row in the first (and only) table and
# apply a yellow highlight across the whole row (paragraph marks + run text),
# matching Word's "highlight" formatting applied via the UI.
$table = $d.Tables.Item(1)
$targetRow = $null
for ($i = 1; $i -le $table.Rows.Count; $i++) {
    $cellText = $table.Rows.Item($i).Cells.Item(1).Range.Text
    if ($cellText -like "nDaysMinBelowZero_Apr*") {
        $targetRow = $table.Rows.Item($i)
        break
    }
}

if ($targetRow -ne $null) {
    $targetRow.Range.Font.HighlightColorIndex = 7
    Write-Output "Highlighted row with nDaysMinBelowZero_Apr"
} else {
    Write-Output "Row not found"
}
